$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Summary": update aggregate stats now that trade #39 (MarketMaking)
# closed and a new trade #72 (MarketMaking) opened.
# ---------------------------------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B4").Value = 0.08     # Total P&L $
$wsSummary.Range("B5").Value = 0.04     # Total P&L %
$wsSummary.Range("B6").Value = 39       # Total Trades
$wsSummary.Range("B8").Value = 17       # Losing Trades
$wsSummary.Range("B9").Value = 41.03    # Win Rate %

# ---------------------------------------------------------------------------
# Sheet "Strategy Status": update the MarketMaking strategy row (row 5).
# ---------------------------------------------------------------------------
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("D5").Value = 6         # Trades
$wsStatus.Range("E5").Value = -0.03     # P&L $
$wsStatus.Range("G5").Value = 33.33     # Win Rate %

# ---------------------------------------------------------------------------
# Sheet "All Trades": close trade #39 (row 40) and append the new trade #72
# (row 73) that was opened right after.
# ---------------------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("All Trades")

# Close out existing trade #39 in row 40
$wsAll.Range("G40").Value = 0.052118
$wsAll.Range("H40").Value = "CLOSED"
$wsAll.Range("I40").Value = -13.1371
$wsAll.Range("J40").Value = -0.01
$wsAll.Range("K40").Value = 100.3
$wsAll.Range("L40").Value = "early_exit"
$wsAll.Range("M40").Value = 0.13

# Append new trade #72 as row 73
$wsAll.Cells.Item(73, 1).Value = 72
$wsAll.Cells.Item(73, 2).NumberFormat = "@"
$wsAll.Cells.Item(73, 2).Value = "2026-02-17"
$wsAll.Cells.Item(73, 3).NumberFormat = "@"
$wsAll.Cells.Item(73, 3).Value = "20:48:58"
$wsAll.Cells.Item(73, 4).Value = "MarketMaking"
$wsAll.Cells.Item(73, 5).Value = "UP"
$wsAll.Cells.Item(73, 6).Value = 0.06
$wsAll.Cells.Item(73, 8).Value = "OPEN"
$wsAll.Cells.Item(73, 9).Value = 0
$wsAll.Cells.Item(73, 10).Value = 0
$wsAll.Cells.Item(73, 11).Value = 100.3038583996649
$wsAll.Cells.Item(73, 13).Value = 0
$wsAll.Cells.Item(73, 14).Value = 0
$wsAll.Cells.Item(73, 15).Value = 0
$wsAll.Cells.Item(73, 16).Value = 0.6
$wsAll.Cells.Item(73, 17).Value = "Normal spread capture: 19600 bps"

# ---------------------------------------------------------------------------
# Sheet "MarketMaking": close trade #39 (row 7) and append the new trade #72
# (row 40) that was opened right after.
# ---------------------------------------------------------------------------
$wsMM = $wb.Worksheets.Item("MarketMaking")

# Close out existing trade #39 in row 7
$wsMM.Range("G7").Value = 0.052118
$wsMM.Range("H7").Value = "CLOSED"
$wsMM.Range("I7").Value = -13.1371
$wsMM.Range("J7").Value = -0.01
$wsMM.Range("K7").Value = 100.3
$wsMM.Range("P7").Value = "early_exit"
$wsMM.Range("Q7").Value = 0.13

# Append new trade #72 as row 40
$wsMM.Cells.Item(40, 1).Value = 72
$wsMM.Cells.Item(40, 2).NumberFormat = "@"
$wsMM.Cells.Item(40, 2).Value = "2026-02-17"
$wsMM.Cells.Item(40, 3).NumberFormat = "@"
$wsMM.Cells.Item(40, 3).Value = "20:48:58"
$wsMM.Cells.Item(40, 4).Value = "MarketMaking"
$wsMM.Cells.Item(40, 5).Value = "UP"
$wsMM.Cells.Item(40, 6).Value = 0.06
$wsMM.Cells.Item(40, 8).Value = "OPEN"
$wsMM.Cells.Item(40, 9).Value = 0
$wsMM.Cells.Item(40, 10).Value = 0
$wsMM.Cells.Item(40, 11).Value = 100.3038583996649
$wsMM.Cells.Item(40, 12).Value = 0
$wsMM.Cells.Item(40, 13).Value = 0
$wsMM.Cells.Item(40, 14).Value = 0.6
$wsMM.Cells.Item(40, 15).Value = "Normal spread capture: 19600 bps"
$wsMM.Cells.Item(40, 17).Value = 0
